$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 42, pushing the existing rows 42 and 43
# (and their data) down to rows 43 and 44 respectively.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record.
$ws.Cells.Item(42, 1).Value = 5
$ws.Cells.Item(42, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(42, 3).Value = "Maule"
$ws.Cells.Item(42, 4).Value = 44491
$ws.Cells.Item(42, 5).Value = 7
$ws.Cells.Item(42, 6).Value = 100112026
$ws.Cells.Item(42, 7).Value = "Haba"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 400
$ws.Cells.Item(42, 11).Value = 8000
$ws.Cells.Item(42, 12).Value = 8000
$ws.Cells.Item(42, 13).Value = 8000
$ws.Cells.Item(42, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(42, 15).Value = "Región del Maule"
$ws.Cells.Item(42, 16).Value = 320
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
